$wb = $excel.ActiveWorkbook

# Sheet references (by position, matching workbook.xml sheet order / rIds)
$wsHtec       = $wb.Worksheets.Item(1)   # HTEC(thousands)
$wsSecThou    = $wb.Worksheets.Item(2)   # <Secondary (thousands)
$wsTerThouTyp = $wb.Worksheets.Item(3)   # <Tertirary (thousands)
$wsTerThou    = $wb.Worksheets.Item(4)   # Tertiary (thousands)
$wsSecPct     = $wb.Worksheets.Item(5)   # <Secondary (%)  -> <Secondary(%)
$wsTerPctTyp  = $wb.Worksheets.Item(6)   # <Tertirary(%)   -> <Tertiary(%)
$wsTerPct     = $wb.Worksheets.Item(7)   # Tertiary(%)

# Rename the two sheets whose tab names had typos/spacing fixed
$wsSecPct.Name = "<Secondary(%)"
$wsTerPctTyp.Name = "<Tertiary(%)"

# Row 10, column A: on HTEC(thousands), <Secondary(%), <Tertiary(%), Tertiary(%)
# the long-form "Germany (until 1990 former territory of the FRG)" label is
# replaced with the shorter "Germany". (Written before the "Country" header
# edits below so the new shared string for "Germany" is appended first.)
$wsHtec.Range("A10").Value = "Germany"
$wsSecPct.Range("A10").Value = "Germany"
$wsTerPctTyp.Range("A10").Value = "Germany"
$wsTerPct.Range("A10").Value = "Germany"

# Row 1, column A: "GEO/TIME" -> "Country" on every sheet. Doing this on all sheets
# removes the last reference to the old shared string so it's dropped entirely,
# shifting every other shared-string index down by one.
$wsHtec.Range("A1").Value = "Country"
$wsSecThou.Range("A1").Value = "Country"
$wsTerThouTyp.Range("A1").Value = "Country"
$wsTerThou.Range("A1").Value = "Country"
$wsSecPct.Range("A1").Value = "Country"
$wsTerPctTyp.Range("A1").Value = "Country"
$wsTerPct.Range("A1").Value = "Country"

# Update the saved selection/scroll position on each sheet to match.
$wsHtec.Range("A2").Select()
$wsSecThou.Range("A1").Select()
$wsTerThouTyp.Range("A2").Select()
$wsTerThou.Range("L10").Select()
$wsSecPct.Range("L11").Select()
$wsTerPctTyp.Range("A2").Select()
$wsTerPct.Range("A2").Select()

# Restore the tab that should remain active/selected.
$wsTerPct.Activate()
